$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "French" test case in A3, pushing the existing "English" value down to A4
$ws.Range("A3").Value = "French"
$ws.Range("A4").Value = "English"

# Update the selection to the newly added row
$ws.Range("A4").Select()
